# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E2) and
# "Correspond Handback DateTime" (H2) values on the zh-cn and de-de
# worksheets to reflect the latest report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 22:58:37"
$wsZhCn.Range("H2").Value = "2016-03-20 22:58:57"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 22:58:40"
$wsDeDe.Range("H2").Value = "2016-03-20 22:59:03"
